$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "to meet; to see (a person) (person に)"
    3  = "there is... (place に thing が)"
    4  = "to buy (～を)"
    5  = "to write (person に thing を)"
    6  = "to take (a picture) (～を)"
    7  = "to wait (～を)"
    8  = "to understand (～が)"
    9  = "(a person) is in...; stays at... (place に person が)"
    40 = "right (～の)"
    41 = "left (～の)"
    42 = "front (～の)"
    43 = "back (～の)"
    44 = "inside (～の)"
    45 = "on (～の)"
    46 = "under (～の)"
    47 = "near; nearby (～の)"
    48 = "next (～の)"
    49 = "between (A と B の)"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
